$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Holly added "S.GISH" to the harvester column (column B) for every data row.
$ws.Range("B2:B13").Value = "S.GISH"

# Match the author's final selection: column B selected, active cell B1.
$ws.Columns("B:B").Select()
